# Adds the "I0" and "IF" columns (I and J) to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (match formatting of the existing header row, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for I (I0) and J (IF) columns, for rows 2 through 63.
$iValues = @(6,7,6,8,8,8,8,9,9,9,9,8,9,9,9,9,9,9,9,9,9,8,7,7,8,8,9,8,11,9,11,9,9,9,8,8,8,8,9,8,6,9,8,8,8,9,9,9,9,8,8,8,8,8,8,8,8,6,6,5,5,4)
$jValues = @(6,7,6,8,8,8,8,9,9,9,9,8,9,9,9,9,9,9,9,9,9,8,8,7,8,8,9,8,11,9,11,9,9,9,8,8,8,8,9,8,6,9,8,8,8,9,9,9,9,8,8,8,8,8,8,8,8,6,6,5,5,4)

for ($r = 2; $r -le 63; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
